$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Cell value edits (impute / un-impute individual missing values)
$ws.Range("E3").Value = -5.7
$ws.Range("E5").ClearContents()
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").ClearContents()

# 2) Fill in the value that used to be missing for SC 193 (currently row 34)
$ws.Range("E34").Value = -6.4

# 3) Remove the two rows that are no longer present in the target data
#    (row 26 "RM 232" and, after that shift, the row that held "SC 92").
#    Delete from bottom-most row first so row numbers of earlier rows
#    stay valid while we work.
$ws.Rows.Item(28).EntireRow.Delete()
$ws.Rows.Item(26).EntireRow.Delete()
